# "Generate Report for Handoff"
#
# A new handoff of b03653a5-8b88-4c4c-af9c-82c32521ab78.md was generated,
# so the localization-status report's handoff timestamps for that file are
# refreshed on all three sheets:
#   - Overview!G6      "Latest HO Xliff Generate Date"
#   - zh-cn!H6         "Latest Handoff Datetime"
#   - de-de!H6         "Latest Handoff Datetime"

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G6").Value = "2016-08-12 08:50:10"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H6").Value = "2016-08-12 08:49:57"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H6").Value = "2016-08-12 08:50:10"
